# Update dashboards - 2025-11-15
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aguilar Prototype")

# --- Row 28 ---
$ws.Range("F28").Value = 0.0292251268148207
$ws.Range("G28").Value = -0.02746655290430811

# --- Row 29 ---
$ws.Range("G29").Value = 0.03397968857635882
$ws.Range("N29").Value = 45975
$ws.Range("Q29").Value = 2.18
$ws.Range("R29").Value = 2.19
$ws.Range("S29").Value = 2.17
$ws.Range("T29").Value = $null
$ws.Range("U29").Value = 2.2

# --- Row 30 ---
$ws.Range("N30").Value = 45975
$ws.Range("R30").Value = 2.28
$ws.Range("S30").Value = 2.27
$ws.Range("T30").Value = $null
$ws.Range("U30").Value = 2.29

# --- Row 47 ---
$ws.Range("N47").Value = 45974
$ws.Range("Q47").Value = 3.88

# --- Row 48 ---
$ws.Range("N48").Value = 45974
$ws.Range("Q48").Value = 3.58
$ws.Range("R48").Value = 3.56
$ws.Range("S48").Value = $null
$ws.Range("T48").Value = 3.58

# --- Row 49 ---
$ws.Range("N49").Value = 45974
$ws.Range("Q49").Value = 3.71
$ws.Range("R49").Value = 3.68
$ws.Range("S49").Value = $null
$ws.Range("T49").Value = 3.72

# --- Row 50 ---
$ws.Range("N50").Value = 45974
$ws.Range("Q50").Value = 4.11
$ws.Range("R50").Value = 4.08
$ws.Range("S50").Value = $null
$ws.Range("T50").Value = 4.13

# --- Row 52 ---
$ws.Range("N52").Value = 45974
$ws.Range("Q52").Value = 5.88
$ws.Range("R52").Value = 5.83
$ws.Range("S52").Value = $null
$ws.Range("T52").Value = 5.86
